# Formatted list with relevant logos - List Format V1 27 02 25
# Adds a new row for the "llama3.2:latest" model, duplicating row 2's
# question/expected-answer/score but with the new model name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3) - same question/expected answer/score as row 2,
# but for the "llama3.2:latest" model.
$ws.Range("A3").Value2 = 3
$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = "llama3.2:latest"
$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2

# Match row 2's formatting (A2 carries the bordered/bold "index" style).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
